# Added receiver RSSI value to out_regs and database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (previously only had the running-sum formula in A11) gains a new
# register entry: size=1, type="u", sql_type="REAL", field="rssi",
# description="RSSI reading for receiver".
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "u"
$ws.Range("D11").Value = "REAL"
$ws.Range("E11").Value = "rssi"
$ws.Range("F11").Value = "RSSI reading for receiver"

# Move the active cell/selection down one row, matching the author's
# recorded selection after making the edit.
$ws.Range("D12").Select()
